# Remove trailing spaces from student surnames in column A
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CHERN"
$ws.Range("A3").Value = "KOH"
$ws.Range("A4").Value = "BRANDON"
$ws.Range("A5").Value = "CALVIN"
$ws.Range("A10").Value = "LEE"
$ws.Range("A11").Value = "LIU"

# Update the active selection to match the post-edit cursor position
$ws.Range("D9").Select()
